# Auto-generated Excel COM-interop script
# Updates cached market-price / profit figures on the Excalibur_Profits sheets
# (ALC, ARM, BSM, CRP, LTW, WVR) to match the latest scheduled-runner pull.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 466.66666
$ws.Range("J2").Value = 800
$ws.Range("L2").Value = 800
$ws.Range("N2").Value = -1026
$ws.Range("H15").Value = 2695.4285
$ws.Range("I15").Value = 2695.4285
$ws.Range("K15").Value = 8086.2855
$ws.Range("M15").Value = -7917.2855
$ws.Range("H62").Value = 14033
$ws.Range("I62").Value = 12015.944
$ws.Range("J62").Value = 20084.166
$ws.Range("K62").Value = 12015.944
$ws.Range("L62").Value = 20084.166
$ws.Range("M62").Value = -11391.944
$ws.Range("N62").Value = -21332.166
$ws.Range("H65").Value = 14033
$ws.Range("I65").Value = 12015.944
$ws.Range("J65").Value = 20084.166
$ws.Range("K65").Value = 60079.72
$ws.Range("L65").Value = 100420.83
$ws.Range("M65").Value = -56959.72
$ws.Range("N65").Value = -106660.83
$ws.Range("H92").Value = 1335.9445
$ws.Range("I92").Value = 407.2857
$ws.Range("K92").Value = 407.2857
$ws.Range("M92").Value = 840.7143
$ws.Range("H99").Value = 1189.1818
$ws.Range("I99").Value = 317.6
$ws.Range("J99").Value = 1915.5
$ws.Range("K99").Value = 952.8000000000001
$ws.Range("L99").Value = 5746.5
$ws.Range("M99").Value = 545.1999999999999
$ws.Range("N99").Value = -8742.5
$ws.Range("H101").Value = 370.5
$ws.Range("I101").Value = 370.5
$ws.Range("K101").Value = 1111.5
$ws.Range("M101").Value = 510.5
$ws.Range("H116").Value = 231800.44
$ws.Range("I116").Value = 8416.666999999999
$ws.Range("K116").Value = 8416.666999999999
$ws.Range("M116").Value = -4974.666999999999
$ws.Range("H118").Value = 3612.8572
$ws.Range("I118").Value = 2572.5
$ws.Range("J118").Value = 5000
$ws.Range("K118").Value = 7717.5
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = -6060.5
$ws.Range("N118").Value = -18314
$ws.Range("H127").Value = 1592.5
$ws.Range("I127").Value = 1592.5
$ws.Range("K127").Value = 4777.5
$ws.Range("M127").Value = 182.5
$ws.Range("H129").Value = 1802.6364
$ws.Range("I129").Value = 1590.6
$ws.Range("J129").Value = 2257
$ws.Range("K129").Value = 4771.799999999999
$ws.Range("L129").Value = 6771
$ws.Range("M129").Value = 228.2000000000007
$ws.Range("N129").Value = -16771
$ws.Range("H138").Value = 8108.3687
$ws.Range("I138").Value = 11052.869
$ws.Range("K138").Value = 33158.607
$ws.Range("M138").Value = -28018.607

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3027.516
$ws.Range("I2").Value = 2483.875
$ws.Range("J2").Value = 4891.4287
$ws.Range("K2").Value = 2483.875
$ws.Range("L2").Value = 4891.4287
$ws.Range("M2").Value = -2370.875
$ws.Range("N2").Value = -5117.4287
$ws.Range("H32").Value = 1717.8148
$ws.Range("I32").Value = 1591.5769
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 1591.5769
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -1304.5769
$ws.Range("N32").Value = -5574
$ws.Range("H33").Value = 2219.2
$ws.Range("I33").Value = 2021.3334
$ws.Range("J33").Value = 4000
$ws.Range("K33").Value = 2021.3334
$ws.Range("L33").Value = 4000
$ws.Range("M33").Value = -1692.3334
$ws.Range("N33").Value = -4658
$ws.Range("H45").Value = 2816.65
$ws.Range("I45").Value = 2816.65
$ws.Range("K45").Value = 2816.65
$ws.Range("M45").Value = -2439.65
$ws.Range("H61").Value = 6324.4443
$ws.Range("I61").Value = 5986.8335
$ws.Range("K61").Value = 5986.8335
$ws.Range("M61").Value = -5774.8335
$ws.Range("H74").Value = 3263.647
$ws.Range("I74").Value = 1597.4166
$ws.Range("J74").Value = 7262.6
$ws.Range("K74").Value = 1597.4166
$ws.Range("L74").Value = 7262.6
$ws.Range("M74").Value = -723.4166
$ws.Range("N74").Value = -9010.6
$ws.Range("H77").Value = 3263.647
$ws.Range("I77").Value = 1597.4166
$ws.Range("J77").Value = 7262.6
$ws.Range("K77").Value = 7987.083000000001
$ws.Range("L77").Value = 36313
$ws.Range("M77").Value = -3619.083000000001
$ws.Range("N77").Value = -45049
$ws.Range("H97").Value = 661.5454999999999
$ws.Range("I97").Value = 661.5454999999999
$ws.Range("K97").Value = 661.5454999999999
$ws.Range("M97").Value = -165.5454999999999
$ws.Range("H116").Value = 3027.516
$ws.Range("I116").Value = 2483.875
$ws.Range("J116").Value = 4891.4287
$ws.Range("K116").Value = 2483.875
$ws.Range("L116").Value = 4891.4287
$ws.Range("M116").Value = -189.875
$ws.Range("N116").Value = -9479.4287
$ws.Range("H122").Value = 3010.353
$ws.Range("I122").Value = 1958.8
$ws.Range("K122").Value = 5876.4
$ws.Range("M122").Value = -3426.4
$ws.Range("H136").Value = 6324.4443
$ws.Range("I136").Value = 5986.8335
$ws.Range("K136").Value = 17960.5005
$ws.Range("M136").Value = -15410.5005

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3027.516
$ws.Range("I3").Value = 2483.875
$ws.Range("J3").Value = 4891.4287
$ws.Range("K3").Value = 2483.875
$ws.Range("L3").Value = 4891.4287
$ws.Range("M3").Value = -2369.875
$ws.Range("N3").Value = -5119.4287
$ws.Range("H20").Value = 1159.1765
$ws.Range("I20").Value = 896.5833
$ws.Range("J20").Value = 1789.4
$ws.Range("K20").Value = 896.5833
$ws.Range("L20").Value = 1789.4
$ws.Range("M20").Value = -649.5833
$ws.Range("N20").Value = -2283.4
$ws.Range("H37").Value = 4000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 4000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 4000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -4274
$ws.Range("H86").Value = 1693.04
$ws.Range("I86").Value = 1538.5264
$ws.Range("K86").Value = 1538.5264
$ws.Range("M86").Value = -415.5264
$ws.Range("H89").Value = 1693.04
$ws.Range("I89").Value = 1538.5264
$ws.Range("K89").Value = 7692.632
$ws.Range("M89").Value = -2076.632
$ws.Range("H99").Value = 11752.615
$ws.Range("I99").Value = 4543.1816
$ws.Range("J99").Value = 51404.5
$ws.Range("K99").Value = 4543.1816
$ws.Range("L99").Value = 51404.5
$ws.Range("M99").Value = -3045.1816
$ws.Range("N99").Value = -54400.5
$ws.Range("H132").Value = 119999.5
$ws.Range("J132").Value = 119999.5
$ws.Range("L132").Value = 119999.5
$ws.Range("N132").Value = -130119.5
$ws.Range("H134").Value = 3651.7817
$ws.Range("I134").Value = 1483.0526
$ws.Range("K134").Value = 4449.1578
$ws.Range("M134").Value = -1914.1578

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 12500
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 8526.5
$ws.Range("J36").Value = 10053
$ws.Range("L36").Value = 10053
$ws.Range("N36").Value = -10829
$ws.Range("H40").Value = 8526.5
$ws.Range("J40").Value = 10053
$ws.Range("L40").Value = 10053
$ws.Range("N40").Value = -10373
$ws.Range("H122").Value = 1896.35
$ws.Range("J122").Value = 2659.3
$ws.Range("L122").Value = 7977.900000000001
$ws.Range("N122").Value = -12877.9

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 589.3125
$ws.Range("I16").Value = 528.1
$ws.Range("J16").Value = 691.3333
$ws.Range("K16").Value = 528.1
$ws.Range("L16").Value = 691.3333
$ws.Range("M16").Value = -358.1
$ws.Range("N16").Value = -1031.3333
$ws.Range("H61").Value = 1419.2858
$ws.Range("I61").Value = 1320.1111
$ws.Range("K61").Value = 1320.1111
$ws.Range("M61").Value = -1118.1111
$ws.Range("H113").Value = 1419.2858
$ws.Range("I113").Value = 1320.1111
$ws.Range("K113").Value = 1320.1111
$ws.Range("M113").Value = 849.8888999999999
$ws.Range("H122").Value = 4885.6665
$ws.Range("I122").Value = 4287.4375
$ws.Range("K122").Value = 12862.3125
$ws.Range("M122").Value = -10412.3125
$ws.Range("H132").Value = 8064.7295
$ws.Range("I132").Value = 8239.571
$ws.Range("K132").Value = 24718.713
$ws.Range("M132").Value = -22188.713
$ws.Range("H136").Value = 4606.2607
$ws.Range("I136").Value = 1617.3636
$ws.Range("K136").Value = 4852.0908
$ws.Range("M136").Value = -2302.0908

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H132").Value = 12647492
$ws.Range("I132").Value = 17462392
$ws.Range("K132").Value = 52387176
$ws.Range("M132").Value = -52384646
